$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row of data for 25/6/2025 into row 44 (text date, like rows 37-43)
$ws.Range("D44").Value = "25/6/2025"
$ws.Range("E44").Value = 297
$ws.Range("F44").Value = 629
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 1012
$ws.Range("J44").Value = "N/A"

# Update the view's selection to match the saved workbook state
$ws.Activate()
$ws.Range("J45").Select()
